$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for rows 2-5. Column D (pool id) stays blank/untouched,
# and columns J-S remain blank as in the source workbook.
# Columns: A, B, C, E, F, G, H, I
$data = @(
    @("sam-1", "sub-1", "Not Defined", "Experimental", "DCE-MRI Contrast Image sam-1", "Breast", "Not Defined", "Not Defined"),
    @("sam-2", "sub-1", "Not Defined", "Experimental", "DCE-MRI Contrast Image sam-2", "Breast", "Not Defined", "Not Defined"),
    @("sam-1", "sub-2", "Not Defined", "Experimental", "DCE-MRI Contrast Image sam-1", "Breast", "Not Defined", "Not Defined"),
    @("sam-2", "sub-2", "Not Defined", "Experimental", "DCE-MRI Contrast Image sam-2", "Breast", "Not Defined", "Not Defined")
)
$columns = @(1, 2, 3, 5, 6, 7, 8, 9)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $col = $columns[$c]
        $ws.Cells.Item($row, $col).Value = $rowValues[$c]
    }
}

# Remove the old trailing rows (formerly rows 6-8: sam-1/sub-2, sam-2/sub-2, sam-3/sub-2)
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(6).Delete()
